# BUG02 Debugging Log.docx edit:
#   - Remove the stray _GoBack bookmark from the title paragraph.
#   - Add a new "Changes to Base Code" section (one Heading1 paragraph and
#     two body paragraphs) right after the "Assumptions" list item ("None"),
#     and move the _GoBack bookmark to the end of the new last paragraph.

$d = $word.ActiveDocument

# --- 1. Strip the "_GoBack" bookmark out of the Title paragraph ------------
$titlePara = $d.Paragraphs(1)
$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:pPr><w:pStyle w:val="Title"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>BUG02</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Debug Log</w:t></w:r>' +
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titlePara.Range.InsertXML($titleXml) | Out-Null

# --- 2. Locate the "None" paragraph (the lone Assumptions bullet) ----------
$noneIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "None") {
        $noneIndex = $i
        break
    }
}

# --- 3. Append a new Heading1 paragraph: "Changes to Base Code" ------------
$noneP = $d.Paragraphs($noneIndex)
$noneP.Range.InsertParagraphAfter() | Out-Null
$headingP = $d.Paragraphs($noneIndex + 1)
$headingXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Changes to Base Code</w:t></w:r>' +
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$headingP.Range.InsertXML($headingXml) | Out-Null

# --- 4. Append the first body paragraph -------------------------------------
$headingP2 = $d.Paragraphs($noneIndex + 1)
$headingP2.Range.InsertParagraphAfter() | Out-Null
$bodyP1 = $d.Paragraphs($noneIndex + 2)
$bodyXml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I have created a substantial refactor to Program.cs on commit:</w:t></w:r>' +
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bodyP1.Range.InsertXML($bodyXml1) | Out-Null

# --- 5. Append the second body paragraph, carrying the _GoBack bookmark ----
$bodyP1b = $d.Paragraphs($noneIndex + 2)
$bodyP1b.Range.InsertParagraphAfter() | Out-Null
$bodyP2 = $d.Paragraphs($noneIndex + 3)
$bodyXml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>This is to separate the structure of the 100 game play, the single game play, and a single round.  This will allow us to tunnel in and test at each level where we need to.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bodyP2.Range.InsertXML($bodyXml2) | Out-Null
